$d = $word.ActiveDocument

# --- Create the three new character styles ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022: Datumi kampanje..." run (4 occurrences) ---

$range = $d.Content
$range.Find.ClearFormatting()
while ($range.Find.Execute("2022: Datumi kampanje za opazovanje ozvezdje Lev: 14.-23. april, 14.-23. maj", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
}

# --- Apply GaNParagraph to the "Sodelujete v svetovni aktivnosti..." run ---

$range2 = $d.Content
$range2.Find.ClearFormatting()
if ($range2.Find.Execute("Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega ozvezdje Lev na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Jenik Hollan, CzechGlobe ..." run ---

$range3 = $d.Content
$range3.Find.ClearFormatting()
if ($range3.Find.Execute("Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range3.Style = "GaNLinks"
}
